# Updates the cryptos price table (Price / Volume(1h) columns) to the
# latest scraped snapshot, mirroring the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.769.92'
$ws.Range("E2").Value = '  +1.81%  '
$ws.Range("D3").Value = '3.286.77'
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.14'
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.49'
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.134'
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("E11").Value = '  +1.66%  '
$ws.Range("D12").Value = '3.862.09'
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("E14").Value = '  +2.34%  '
$ws.Range("D15").Value = '68.785.22'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000172'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("D17").Value = '3.271.76'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.58'
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '394.54'
$ws.Range("E20").Value = '  +4.52%  '
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.66'
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.516'
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("E25").Value = '  +1.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.188'
$ws.Range("E26").Value = '  +3.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.65'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.04'
$ws.Range("E31").Value = '  +1.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.15'
$ws.Range("E32").Value = '  +3.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.30'
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  +1.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.88'
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("E37").Value = '  +1.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.835'
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("E39").Value = '  +2.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.18'
$ws.Range("E40").Value = '  -2.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.55'
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("E42").Value = '  -2.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.53'
$ws.Range("E43").Value = '  +2.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '343.33'
$ws.Range("E45").Value = '  -5.44%  '
$ws.Range("D46").Value = '2.610.32'
$ws.Range("E46").Value = '  -4.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.91'
$ws.Range("E47").Value = '  -2.22%  '
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.01'
$ws.Range("E49").Value = '  +3.53%  '
$ws.Range("E50").Value = '  +2.75%  '
$ws.Range("E51").Value = '  -0.10%  '
